$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.004.92"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.957.13"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.65%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "379.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.87"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.49%  "
$ws.Range("E7").Value = "  +0.69%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.585"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.31"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.66%  "
$ws.Range("E11").Value = "  -0.29%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0850"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "12.75"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +78.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.43"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.34%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.423.80"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.54%  "
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.79"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.46%  "
$ws.Range("B17").Value = "Polygon"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.02"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.44%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.963.55"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "50.958.99"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("E20").Value = "  -2.97%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.41"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.10%  "
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("B23").Value = "PancakeSwap"
$ws.Range("C23").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +13.05%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.56"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.54%  "
$ws.Range("B25").Value = "BitcoinCash"
$ws.Range("C25").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "266.63"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.82%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.10"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.99%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "25.72"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.35%  "
$ws.Range("E29").Value = "  -2.71%  "
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.02"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -8.76%  "
$ws.Range("E31").Value = "  -4.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.29"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.02%  "
$ws.Range("E33").Value = "  +0.12%  "
$ws.Range("B34").Value = "Toncoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.05"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.33%  "
$ws.Range("B35").Value = "InjectiveProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "33.86"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.11%  "
$ws.Range("E36").Value = "  -4.89%  "
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.14"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.25%  "
$ws.Range("E39").Value = "  +1.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.59"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.16%  "
$ws.Range("E41").Value = "  +2.37%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.50"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.22%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "119.92"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.79%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.57"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +10.49%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.43"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.02"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.39%  "
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.009.38"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.19%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.31"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.66%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.259"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.85%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0320"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -7.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.24%  "
